# Journal de bord Julien - add the "Samedi 1 novembre 2014" entry after
# the existing "Vendredi 31 octobre 2011" entry.
#
# The new paragraph is built as a WordprocessingML fragment (so that we can
# reproduce the exact run layout - bold date label, plain body runs, and the
# spell-check proofErr wrappers around the camel-case class names CLabel /
# CLabelLeftRight / CButton / CControl) and inserted via Range.InsertXML at
# the very end of the document, right where the old "Vendredi" paragraph
# (and its trailing _GoBack bookmark) used to end.
#
# The _GoBack bookmark always sits on the last edited spot in the document,
# so it is removed from its old position and re-created (embedded directly
# in the inserted fragment) at the end of the brand-new paragraph, matching
# the target layout.

$d = $word.ActiveDocument

if ($d.Bookmarks.Exists("_GoBack")) {
    $gb = $d.Bookmarks("_GoBack")
    $gb.Delete()
}

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Samedi 1 novembre 2014 : </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">J’ai créé la fenêtre du programme ainsi que les menus principal et nouvelle partie. J’ai terminé la création du menu principal et tout s’affiche bien. Les classes </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>CLabel</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>CLabelLeftRight</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> et </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>CButton</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> hérite maintenant de la classe </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>CControl</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$insertionPoint.InsertXML($xml)
